$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared string "testing_route" -> "acma_check" wherever it is used (B2, B3)
$ws.Range("B2").Value = "acma_check"
$ws.Range("B3").Value = "acma_check"

# Update the active cell selection from F16 to D6
$ws.Range("D6").Select()
